$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - apostrophe-prefixed to force text and preserve exact formatting
$priceUpdates = @{
    'D2' = '78.816.17'
    'D3' = '3.187.11'
    'D5' = '205.41'
    'D6' = '629.74'
    'D10' = '3.187.60'
    'D11' = '0.585'
    'D12' = '0.164'
    'D13' = '5.44'
    'D14' = '3.774.12'
    'D15' = '0.0000225'
    'D16' = '31.43'
    'D17' = '78.755.75'
    'D18' = '3.191.86'
    'D19' = '14.47'
    'D20' = '9.24'
    'D21' = '427.75'
    'D22' = '2.87'
    'D23' = '4.91'
    'D25' = '3.349.11'
    'D27' = '11.01'
    'D28' = '75.72'
    'D31' = '1.00'
    'D32' = '8.85'
    'D34' = '514.60'
    'D37' = '22.96'
    'D39' = '0.999'
    'D41' = '163.90'
    'D43' = '192.22'
    'D46' = '0.804'
    'D48' = '1.30'
    'D49' = '42.65'
}

# Volume(1h) (column E) updates
$volumeUpdates = @{
    'E2' = '  +3.11%  '
    'E3' = '  +5.39%  '
    'E4' = '  +0.00%  '
    'E5' = '  +2.62%  '
    'E6' = '  +0.17%  '
    'E7' = '  +0.04%  '
    'E8' = '  +8.95%  '
    'E9' = '  +5.36%  '
    'E10' = '  +5.42%  '
    'E11' = '  +34.24%  '
    'E12' = '  +2.60%  '
    'E13' = '  +5.93%  '
    'E14' = '  +5.41%  '
    'E15' = '  +17.34%  '
    'E16' = '  +7.68%  '
    'E18' = '  +4.98%  '
    'E19' = '  +7.57%  '
    'E20' = '  +2.11%  '
    'E21' = '  +14.51%  '
    'E22' = '  +26.65%  '
    'E23' = '  +12.61%  '
    'E25' = '  +5.11%  '
    'E26' = '  +8.02%  '
    'E27' = '  +11.38%  '
    'E28' = '  +3.69%  '
    'E29' = '  +0.19%  '
    'E30' = '  +5.32%  '
    'E31' = '  +0.19%  '
    'E32' = '  +6.63%  '
    'E33' = '  +4.85%  '
    'E34' = '  +1.55%  '
    'E35' = '  +2.00%  '
    'E36' = '  +21.50%  '
    'E37' = '  +10.85%  '
    'E38' = '  +19.07%  '
    'E39' = '  -0.01%  '
    'E40' = '  +3.42%  '
    'E41' = '  -0.11%  '
    'E42' = '  -0.12%  '
    'E43' = '  -0.27%  '
    'E44' = '  -0.31%  '
    'E45' = '  +6.28%  '
    'E46' = '  +12.12%  '
    'E47' = '  +6.62%  '
    'E48' = '  +3.88%  '
    'E49' = '  +0.70%  '
    'E50' = '  +5.02%  '
    'E51' = '  +3.41%  '
}

foreach ($cell in $priceUpdates.Keys) {
    $ws.Range($cell).Value = "'" + $priceUpdates[$cell]
}

foreach ($cell in $volumeUpdates.Keys) {
    $ws.Range($cell).Value = $volumeUpdates[$cell]
}
